# Apply the edits described in the commit "adding template without tests".
#
# Real content changes:
#  1. On the "Exclude" sheet, the "size" parameter's value list (column G,
#     rows 2/5/9/12) gains two more bucket values: ...,10,3 -> ...,10,3,2000,2400
#  2. The "Exclude" sheet becomes the active / selected sheet (with cell
#     G12 selected) instead of "BayCountKPI".
#  3. Column G on "Exclude" is widened (it holds the longer value list now).
#  4. Page setup across sheets is unified to paper size 9 (A4).
#  5. Iterative-calculation max change (iterateDelta) is tightened.
#  6. Workbook tab-bar ratio is adjusted.

$wb = $excel.ActiveWorkbook

$wsExclude     = $wb.Worksheets.Item("Exclude")
$wsInclude     = $wb.Worksheets.Item("Include")
$wsBayCountKPI = $wb.Worksheets.Item("BayCountKPI")

# 1. Extend the "size" value list wherever it appears on the Exclude sheet.
$newSizeList = "2,4,2.25,2.4,10,3,2000,2400"
$wsExclude.Range("G2").Value = $newSizeList
$wsExclude.Range("G5").Value = $newSizeList
$wsExclude.Range("G9").Value = $newSizeList
$wsExclude.Range("G12").Value = $newSizeList

# 3. Column G needs to be wider to fit the longer list of values.
$wsExclude.Columns.Item(7).ColumnWidth = 20.5

# 4. Normalise the printed paper size to A4 (9) on the sheets that still
#    had the old default (1).
$wsExclude.PageSetup.PaperSize = 9
$wsInclude.PageSetup.PaperSize = 9

# 5. Tighten the iterative-calculation convergence threshold.
try {
    $excel.IterativeCalculation = $true
    $excel.MaxIterations = 100
    $excel.MaxChange = 0.0001
} catch {}

# 6. Tab ratio (tab bar vs horizontal scroll bar split) moves to 50%.
try {
    $excel.ActiveWindow.TabRatio = 0.5
} catch {}

# 2. Make "Exclude" the active sheet with G12 selected, leaving
#    "BayCountKPI" no longer the selected tab.
$wsExclude.Activate() | Out-Null
$wsExclude.Range("G12").Select() | Out-Null
